# Student import template rework:
#  - Collapse the two-sheet "Instructions" + "Template" workbook into a
#    single "Students" sheet that is ready to fill in and upload.
#  - Fix up the sample data (passwords / gradeLevel) and drop the third
#    sample row (Michael Johnson), leaving just John + Jane as examples.
#  - Give the columns sensible, readable widths.

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Template")
$instructions = $wb.Worksheets.Item("Instructions")

# Move the template's data onto the sheet that will survive (keeps the
# workbook down to a single sheet instead of just deleting+renaming, which
# is also closer to what a user re-purposing the first tab would do).
$instructions.Cells.Clear()
$template.Range("A1:H4").Copy($instructions.Range("A1"))

# The sample sheet has done its job - drop it.
$template.Delete()

# This is now the one and only sheet - rename it to match its new purpose.
$ws = $instructions
$ws.Name = "Students"

# --- Fix up the sample rows --------------------------------------------------
$ws.Range("E2").Value = "Password123"
$ws.Range("F2").Value = "'9"

$ws.Range("E3").Value = "Password123"
$ws.Range("F3").Value = "'9"

# Remove the third sample row (Michael Johnson) so only John + Jane remain.
$ws.Rows(4).Delete()

# --- Column widths ------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15
$ws.Columns.Item(2).ColumnWidth = 15
$ws.Columns.Item(3).ColumnWidth = 20
$ws.Columns.Item(4).ColumnWidth = 25
$ws.Columns.Item(5).ColumnWidth = 15
$ws.Columns.Item(6).ColumnWidth = 10
$ws.Columns.Item(7).ColumnWidth = 10
$ws.Columns.Item(8).ColumnWidth = 8
